# Quest.xlsx update: fill in the missing "Reward Item Id" (column H) values
# for the quest rows that didn't have one yet, and move the active
# selection to I7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 3, 4, 6, 7 and 10 are missing an H value; set them to 0.
$rowsToFill = @(3, 4, 6, 7, 10)
foreach ($r in $rowsToFill) {
    $ws.Cells.Item($r, 8).Value = 0
}

# Move the current selection to I7.
$ws.Range("I7").Select()
